# Updates cryptos list price/volume data (and the two row swaps for ranks 10/11 and 44/45)
# per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.134.83"
$ws.Range("E2").Value = "  -1.62%  "
# Row 3
$ws.Range("D3").Value = "1.653.06"
$ws.Range("E3").Value = "  -1.96%  "
# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.45%  "
# Row 5
$ws.Range("D5").Value = "'218.38"
$ws.Range("E5").Value = "  +0.37%  "
# Row 6
$ws.Range("D6").Value = "'0.5205"
$ws.Range("E6").Value = "  -2.66%  "
# Row 7
$ws.Range("E7").Value = "  +0.39%  "
# Row 8
$ws.Range("D8").Value = "'0.2668"
$ws.Range("E8").Value = "  -0.76%  "
# Row 9
$ws.Range("D9").Value = "'0.06326"
$ws.Range("E9").Value = "  -1.77%  "
# Row 10
$ws.Range("D10").Value = "'21.04"
$ws.Range("E10").Value = "  -1.92%  "
# Row 11
$ws.Range("D11").Value = "'0.07738"
$ws.Range("E11").Value = "  -0.54%  "
# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.444"
$ws.Range("E12").Value = "  -1.39%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.653.86"
$ws.Range("E13").Value = "  -1.72%  "
# Row 14
$ws.Range("D14").Value = "1.881.37"
$ws.Range("E14").Value = "  -1.84%  "
# Row 15
$ws.Range("D15").Value = "'0.5457"
$ws.Range("E15").Value = "  -3.14%  "
# Row 16
$ws.Range("D16").Value = "0.0₅8226"
$ws.Range("E16").Value = "  -2.59%  "
# Row 17
$ws.Range("D17").Value = "'64.82"
$ws.Range("E17").Value = "  -1.99%  "
# Row 18
$ws.Range("D18").Value = "26.178.79"
$ws.Range("E18").Value = "  -1.55%  "
# Row 19
$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "  +0.46%  "
# Row 20
$ws.Range("D20").Value = "'4.666"
$ws.Range("E20").Value = "  -3.27%  "
# Row 21
$ws.Range("D21").Value = "'192.62"
$ws.Range("E21").Value = "  -1.00%  "
# Row 22
$ws.Range("E22").Value = "  -2.55%  "
# Row 23
$ws.Range("D23").Value = "'6.091"
$ws.Range("E23").Value = "  -4.91%  "
# Row 24
$ws.Range("D24").Value = "'1.009"
$ws.Range("E24").Value = "  +0.65%  "
# Row 25
$ws.Range("D25").Value = "'137.48"
$ws.Range("E25").Value = "  -4.56%  "
# Row 26
$ws.Range("D26").Value = "'0.1235"
$ws.Range("E26").Value = "  -3.08%  "
# Row 27
$ws.Range("D27").Value = "'7.238"
$ws.Range("E27").Value = "  -3.40%  "
# Row 28
$ws.Range("D28").Value = "'16.11"
$ws.Range("E28").Value = "  -0.71%  "
# Row 29
$ws.Range("D29").Value = "'1.408"
$ws.Range("E29").Value = "  -0.81%  "
# Row 30
$ws.Range("D30").Value = "'0.06035"
$ws.Range("E30").Value = "  -1.37%  "
# Row 31
$ws.Range("D31").Value = "'1.283"
$ws.Range("E31").Value = "  +0.32%  "
# Row 32
$ws.Range("D32").Value = "'3.567"
$ws.Range("E32").Value = "  -1.10%  "
# Row 33
$ws.Range("D33").Value = "'3.337"
$ws.Range("E33").Value = "  -3.76%  "
# Row 34
$ws.Range("D34").Value = "'1.650"
$ws.Range("E34").Value = "  -3.06%  "
# Row 35
$ws.Range("D35").Value = "'0.9802"
$ws.Range("E35").Value = "  -3.88%  "
# Row 36
$ws.Range("D36").Value = "'2.412"
$ws.Range("E36").Value = "  -0.27%  "
# Row 37
$ws.Range("D37").Value = "'2.773"
$ws.Range("E37").Value = "  -0.84%  "
# Row 38
$ws.Range("D38").Value = "'0.5914"
$ws.Range("E38").Value = "  +3.39%  "
# Row 39
$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "  -3.02%  "
# Row 40
$ws.Range("D40").Value = "'5.951"
$ws.Range("E40").Value = "  -0.43%  "
# Row 41
$ws.Range("D41").Value = "'0.8635"
$ws.Range("E41").Value = "  -0.83%  "
# Row 42
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  +0.26%  "
# Row 43
$ws.Range("D43").Value = "1.036.70"
$ws.Range("E43").Value = "  -1.95%  "
# Row 44
$ws.Range("D44").Value = "'99.72"
$ws.Range("E44").Value = "  -0.47%  "
# Row 45
$ws.Range("D45").Value = "1.794.63"
$ws.Range("E45").Value = "  -2.28%  "
# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'57.15"
$ws.Range("E46").Value = "  -0.24%  "
# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈108"
$ws.Range("E47").Value = "  -4.45%  "
# Row 48
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  +0.13%  "
# Row 49
$ws.Range("E49").Value = "  -0.41%  "
# Row 50
$ws.Range("E50").Value = "  -0.53%  "
# Row 51
$ws.Range("D51").Value = "'1.465"
$ws.Range("E51").Value = "  +3.49%  "
